# Perbaikan Klast dan PCA
# Update column U (BBL) values for rows 2-13 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 299.96
    3  = 285.55
    4  = 354.79
    5  = 441.2
    6  = 453.19
    7  = 360.17
    8  = 269.65
    9  = 293.42
    10 = 263.66
    11 = 333.76
    12 = 307.42
    13 = 294.61
}

foreach ($row in $updates.Keys) {
    $ws.Range("U$row").Value = $updates[$row]
}
